# Update Sheets via scheduled runner
# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# (columns H:N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR
# sheets, row-by-row, matching the latest market-board pull.

$wb = $excel.ActiveWorkbook

$changes = @(
    ,@("ALC","H51",4674.5)
    ,@("ALC","I51",2445)
    ,@("ALC","J51",4922.222)
    ,@("ALC","K51",2445)
    ,@("ALC","L51",4922.222)
    ,@("ALC","M51",-1961)
    ,@("ALC","N51",-5890.222)
    ,@("ALC","H80",687.3125)
    ,@("ALC","I80",464.77777)
    ,@("ALC","J80",973.4286)
    ,@("ALC","K80",1394.33331)
    ,@("ALC","L80",2920.2858)
    ,@("ALC","M80",-396.33331)
    ,@("ALC","N80",-4916.2858)
    ,@("ALC","H83",687.3125)
    ,@("ALC","I83",464.77777)
    ,@("ALC","J83",973.4286)
    ,@("ALC","K83",4182.99993)
    ,@("ALC","L83",8760.857399999999)
    ,@("ALC","M83",809.0000700000001)
    ,@("ALC","N83",-18744.8574)
    ,@("ALC","H112",2585.5293)
    ,@("ALC","J112",2585.5293)
    ,@("ALC","L112",7756.5879)
    ,@("ALC","N112",-9972.5879)
    ,@("ALC","H137",22729160)
    ,@("ALC","I137",35716180)
    ,@("ALC","J137",1876.125)
    ,@("ALC","K137",107148540)
    ,@("ALC","L137",5628.375)
    ,@("ALC","M137",-107145990)
    ,@("ALC","N137",-10728.375)
    ,@("ARM","H2",1416.2142)
    ,@("ARM","I2",1332.9231)
    ,@("ARM","K2",1332.9231)
    ,@("ARM","M2",-1219.9231)
    ,@("ARM","H32",3948.879)
    ,@("ARM","I32",3107.1292)
    ,@("ARM","K32",3107.1292)
    ,@("ARM","M32",-2820.1292)
    ,@("ARM","H45",2000.909)
    ,@("ARM","I45",2399.25)
    ,@("ARM","J45",1773.2858)
    ,@("ARM","K45",2399.25)
    ,@("ARM","L45",1773.2858)
    ,@("ARM","M45",-2022.25)
    ,@("ARM","N45",-2527.2858)
    ,@("ARM","H97",495.45)
    ,@("ARM","I97",465.35715)
    ,@("ARM","K97",465.35715)
    ,@("ARM","M97",30.64285000000001)
    ,@("ARM","H110",1271.2)
    ,@("ARM","I110",1144.5714)
    ,@("ARM","J110",1566.6666)
    ,@("ARM","K110",1144.5714)
    ,@("ARM","L110",1566.6666)
    ,@("ARM","M110",900.4286)
    ,@("ARM","N110",-5656.6666)
    ,@("ARM","H116",1416.2142)
    ,@("ARM","I116",1332.9231)
    ,@("ARM","K116",1332.9231)
    ,@("ARM","M116",961.0769)
    ,@("BSM","H3",1416.2142)
    ,@("BSM","I3",1332.9231)
    ,@("BSM","K3",1332.9231)
    ,@("BSM","M3",-1218.9231)
    ,@("BSM","H105",1819.7222)
    ,@("BSM","I105",1740.5)
    ,@("BSM","K105",1740.5)
    ,@("BSM","M105",6.5)
    ,@("CRP","H22",956.3684)
    ,@("CRP","I22",363.18182)
    ,@("CRP","J22",1772)
    ,@("CRP","K22",363.18182)
    ,@("CRP","L22",1772)
    ,@("CRP","M22",-13.18182000000002)
    ,@("CRP","N22",-2472)
    ,@("CRP","H31",1838.5186)
    ,@("CRP","J31",4797.6)
    ,@("CRP","L31",4797.6)
    ,@("CRP","N31",-5387.6)
    ,@("CRP","H34",1838.5186)
    ,@("CRP","J34",4797.6)
    ,@("CRP","L34",4797.6)
    ,@("CRP","N34",-5201.6)
    ,@("CRP","H132",3653)
    ,@("CRP","I132",2979.5)
    ,@("CRP","K132",8938.5)
    ,@("CRP","M132",-6408.5)
    ,@("CRP","H134",1611.8572)
    ,@("CRP","I134",1398.1666)
    ,@("CRP","J134",1772.125)
    ,@("CRP","K134",4194.4998)
    ,@("CRP","L134",5316.375)
    ,@("CRP","M134",-1659.4998)
    ,@("CRP","N134",-10386.375)
    ,@("CUL","H2",11905141)
    ,@("CUL","I2",232.2)
    ,@("CUL","J2",25641574)
    ,@("CUL","K2",1393.2)
    ,@("CUL","L2",153849444)
    ,@("CUL","M2",-1280.2)
    ,@("CUL","N2",-153849670)
    ,@("CUL","H22",1333.3334)
    ,@("CUL","I22",1000)
    ,@("CUL","J22",1500)
    ,@("CUL","K22",3000)
    ,@("CUL","L22",4500)
    ,@("CUL","M22",-2831)
    ,@("CUL","N22",-4838)
    ,@("CUL","H27",1333.3334)
    ,@("CUL","I27",1000)
    ,@("CUL","J27",1500)
    ,@("CUL","K27",3000)
    ,@("CUL","L27",4500)
    ,@("CUL","M27",-2898)
    ,@("CUL","N27",-4704)
    ,@("CUL","H37",117338.164)
    ,@("CUL","J37",117338.164)
    ,@("CUL","L37",352014.492)
    ,@("CUL","N37",-352238.492)
    ,@("CUL","H39",3371.1428)
    ,@("CUL","J39",4459.6)
    ,@("CUL","L39",13378.8)
    ,@("CUL","N39",-13966.8)
    ,@("CUL","H44",474)
    ,@("CUL","I44",474)
    ,@("CUL","J44",0)
    ,@("CUL","K44",1422)
    ,@("CUL","L44",0)
    ,@("CUL","M44",-1024)
    ,@("CUL","N44",$null)
    ,@("CUL","H131",2198.75)
    ,@("CUL","I131",2465.3333)
    ,@("CUL","K131",7395.999899999999)
    ,@("CUL","M131",-2355.999899999999)
    ,@("GSM","H33",39750.25)
    ,@("GSM","J33",52666.668)
    ,@("GSM","L33",52666.668)
    ,@("GSM","N33",-53170.668)
    ,@("GSM","H70",32429.637)
    ,@("GSM","I70",50533.082)
    ,@("GSM","J70",10705.5)
    ,@("GSM","K70",50533.082)
    ,@("GSM","L70",10705.5)
    ,@("GSM","M70",-50263.082)
    ,@("GSM","H73",32429.637)
    ,@("GSM","I73",50533.082)
    ,@("GSM","J73",10705.5)
    ,@("GSM","K73",50533.082)
    ,@("GSM","L73",10705.5)
    ,@("GSM","M73",-49597.082)
    ,@("GSM","H97",1500)
    ,@("GSM","I97",0)
    ,@("GSM","J97",1500)
    ,@("GSM","K97",0)
    ,@("GSM","L97",1500)
    ,@("GSM","M97",$null)
    ,@("GSM","N97",-2492)
    ,@("GSM","H102",2114.125)
    ,@("GSM","I102",1993.7858)
    ,@("GSM","J102",2956.5)
    ,@("GSM","K102",1993.7858)
    ,@("GSM","L102",2956.5)
    ,@("GSM","M102",-371.7858000000001)
    ,@("GSM","N102",-6200.5)
    ,@("GSM","H121",26998)
    ,@("GSM","J121",26998)
    ,@("GSM","L121",26998)
    ,@("GSM","N121",-30492)
    ,@("GSM","H126",8554.888999999999)
    ,@("GSM","I126",8499.75)
    ,@("GSM","K126",25499.25)
    ,@("GSM","M126",-23029.25)
    ,@("LTW","H16",1670.7693)
    ,@("LTW","I16",1670.7693)
    ,@("LTW","K16",1670.7693)
    ,@("LTW","M16",-1500.7693)
    ,@("LTW","H22",63837.59)
    ,@("LTW","I22",250612.25)
    ,@("LTW","K22",250612.25)
    ,@("LTW","M22",-250317.25)
    ,@("LTW","H27",63837.59)
    ,@("LTW","I27",250612.25)
    ,@("LTW","K27",250612.25)
    ,@("LTW","M27",-250505.25)
    ,@("LTW","H55",1171.5555)
    ,@("LTW","I55",542.4286)
    ,@("LTW","K55",542.4286)
    ,@("LTW","M55",-369.4286)
    ,@("LTW","H122",7812.375)
    ,@("LTW","I122",5749.1665)
    ,@("LTW","J122",14002)
    ,@("LTW","K122",17247.4995)
    ,@("LTW","L122",42006)
    ,@("LTW","M122",-14797.4995)
    ,@("LTW","N122",-46906)
    ,@("LTW","H136",4831.6924)
    ,@("LTW","I136",4155.6665)
    ,@("LTW","J136",5411.143)
    ,@("LTW","K136",12466.9995)
    ,@("LTW","L136",16233.429)
    ,@("LTW","M136",-9916.999500000002)
    ,@("LTW","N136",-21333.429)
    ,@("WVR","H108",0)
    ,@("WVR","J108",0)
    ,@("WVR","L108",0)
    ,@("WVR","N108",$null)
    ,@("WVR","H109",0)
    ,@("WVR","J109",0)
    ,@("WVR","L109",0)
    ,@("WVR","N109",$null)
    ,@("WVR","H122",14517810)
    ,@("WVR","I122",10418580)
    ,@("WVR","K122",31255740)
    ,@("WVR","M122",-31253290)
    ,@("WVR","H132",1979.1666)
    ,@("WVR","I132",1775)
    ,@("WVR","K132",5325)
    ,@("WVR","M132",-2795)
    ,@("WVR","H136",6821.875)
    ,@("WVR","I136",8250.091)
    ,@("WVR","J136",3679.8)
    ,@("WVR","K136",8250.091)
    ,@("WVR","L136",11039.4)
    ,@("WVR","M136",-22200.273)
    ,@("WVR","N136",-16139.4)
)

$wsCache = @{}

foreach ($item in $changes) {
    $sheetName = $item[0]
    $addr      = $item[1]
    $val       = $item[2]

    if (-not $wsCache.ContainsKey($sheetName)) {
        $wsCache[$sheetName] = $wb.Worksheets.Item($sheetName)
    }
    $ws = $wsCache[$sheetName]

    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}
